# Cristofor_Rotsching_AAC2_Harta_particule_elementare.pptx
#
# 1) SlideMaster3 title placeholder: the placeholder's prompt text was
#    authored as one <a:r> run per character ("C","l","i","c","k", " ", ...).
#    Re-set it as a single run containing the full phrase.
# 2) Slide 10, shape "CustomShape 1": drop the third paragraph's run/text
#    "(Functia Dirac delta)" while keeping the (now empty) paragraph.

$p = $ppt.ActivePresentation

# --- 1. Fix the title placeholder prompt text on slide master 3 -----------
$design = $p.Designs.Item(3)
$master = $design.SlideMaster
$titlePlaceholder = $master.Shapes.Item(3)
$titlePlaceholder.TextFrame.TextRange.Text = "Click to edit the title text format"

# --- 2. Remove the "(Functia Dirac delta)" run from slide 10 --------------
$slide = $p.Slides.Item(10)
$shape = $slide.Shapes.Item(1)
$enDash = [char]0x2013
$shape.TextFrame.TextRange.Text = "Procesarea $enDash Semnalului`rSmoothing`r"
